$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.389.40'
$ws.Range('E2').Value = '  -1.49%  '

$ws.Range('D3').Value = '2.369.40'
$ws.Range('E3').Value = '  +0.84%  '

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').Value = "'" + '331.12'
$ws.Range('E5').Value = '  +6.11%  '

$ws.Range('D6').Value = "'" + '100.23'
$ws.Range('E6').Value = '  -6.89%  '

$ws.Range('D7').Value = "'" + '0.637'
$ws.Range('E7').Value = '  -0.29%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = "'" + '0.619'
$ws.Range('E9').Value = '  -0.62%  '

$ws.Range('D10').Value = "'" + '40.27'
$ws.Range('E10').Value = '  -6.74%  '

$ws.Range('E11').Value = '  -2.11%  '

$ws.Range('D12').Value = "'" + '8.47'
$ws.Range('E12').Value = '  -5.27%  '

$ws.Range('D13').Value = "'" + '1.02'
$ws.Range('E13').Value = '  -4.54%  '

$ws.Range('D14').Value = "'" + '0.106'
$ws.Range('E14').Value = '  +0.51%  '

$ws.Range('D15').Value = "'" + '16.30'
$ws.Range('E15').Value = '  -0.51%  '

$ws.Range('D16').Value = '2.723.29'
$ws.Range('E16').Value = '  +0.97%  '

$ws.Range('D17').Value = '2.362.54'
$ws.Range('E17').Value = '  +1.10%  '

$ws.Range('D18').Value = '42.506.69'
$ws.Range('E18').Value = '  -1.07%  '

$ws.Range('D19').Value = "'" + '7.74'
$ws.Range('E19').Value = '  +6.61%  '

$ws.Range('E20').Value = '  -1.50%  '

$ws.Range('B21').Value = 'PancakeSwap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D21').Value = "'" + '3.76'
$ws.Range('E21').Value = '  +9.19%  '

$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').Value = "'" + '75.05'
$ws.Range('E22').Value = '  -0.68%  '

$ws.Range('D23').Value = "'" + '275.59'
$ws.Range('E23').Value = '  +8.42%  '

$ws.Range('E24').Value = '  -8.57%  '

$ws.Range('D25').Value = "'" + '9.71'
$ws.Range('E25').Value = '  +9.12%  '

$ws.Range('E26').Value = '  +0.16%  '

$ws.Range('D27').Value = "'" + '11.49'
$ws.Range('E27').Value = '  -4.35%  '

$ws.Range('D28').Value = "'" + '23.64'
$ws.Range('E28').Value = '  +3.66%  '

$ws.Range('E29').Value = '  -1.33%  '

$ws.Range('D30').Value = "'" + '174.09'
$ws.Range('E30').Value = '  +0.75%  '

$ws.Range('E31').Value = '  -2.14%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'" + '0.0902'
$ws.Range('E32').Value = '  -0.75%  '

$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = "'" + '35.35'
$ws.Range('E33').Value = '  -9.40%  '

$ws.Range('D34').Value = "'" + '6.04'
$ws.Range('E34').Value = '  +3.55%  '

$ws.Range('D35').Value = "'" + '0.134'
$ws.Range('E35').Value = '  +2.26%  '

$ws.Range('D36').Value = "'" + '4.59'
$ws.Range('E36').Value = '  -7.24%  '

$ws.Range('D37').Value = "'" + '0.0360'
$ws.Range('E37').Value = '  -4.34%  '

$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = "'" + '2.93'
$ws.Range('E38').Value = '  +5.63%  '

$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = "'" + '3.87'
$ws.Range('E39').Value = '  -6.02%  '

$ws.Range('E40').Value = '  +1.61%  '

$ws.Range('D41').Value = "'" + '1.53'
$ws.Range('E41').Value = '  +1.06%  '

$ws.Range('D42').Value = "'" + '0.229'
$ws.Range('E42').Value = '  -1.47%  '

$ws.Range('D43').Value = "'" + '69.35'
$ws.Range('E43').Value = '  -3.76%  '

$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'" + '115.81'
$ws.Range('E45').Value = '  +3.33%  '

$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').Value = "'" + '89.23'
$ws.Range('E46').Value = '  +29.40%  '

$ws.Range('D47').Value = "'" + '11.93'
$ws.Range('E47').Value = '  -3.85%  '

$ws.Range('D48').Value = "'" + '5.47'
$ws.Range('E48').Value = '  -3.24%  '

$ws.Range('D49').Value = "'" + '9.10'
$ws.Range('E49').Value = '  -1.18%  '

$ws.Range('D50').Value = '1.593.52'
$ws.Range('E50').Value = '  +7.05%  '

$ws.Range('D51').Value = "'" + '1.27'
$ws.Range('E51').Value = '  -2.16%  '
